# Add two new taxon records ("Liriomyza equiseti" and "Chirosia similata")
# to the "New Taxa" worksheet of the FWSpecies additions workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Taxa")

# --- Row 5: Liriomyza equiseti De Meijere, 1924 ---------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Insect|Species"
$ws.Range("D5").Value = "Liriomyza equiseti De Meijere, 1924"
$ws.Range("C5").Value = "Liriomyza equiseti"
$ws.Range("F5").Value = "accepted"
$ws.Range("G5").Value = 49281
$ws.Range("L5").Value = "<em>Liriomyza</em> <em>equiseti</em> De Meijere, 1924"
$ws.Range("K5").Value = "<em>Liriomyza</em> <em>equiseti</em>"
$ws.Range("M5").Value = "https://www.gbif.org/species/1553334"

# --- Row 6: Chirosia similata (Tiensuu, 1939) ------------------------------
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Insect|Species"
$ws.Range("D6").Value = "Chirosia similata (Tiensuu, 1939)"
$ws.Range("C6").Value = "Chirosia similata"
$ws.Range("F6").Value = "accepted"
$ws.Range("G6").Value = 56467
$ws.Range("L6").Value = "<em>Chirosia</em> <em>similata</em> (Tiensuu, 1939)"
$ws.Range("K6").Value = "<em>Chirosia</em> <em>similata</em>"
$ws.Range("M6").Value = "https://www.gbif.org/species/1575323"

# Match the author's final selection state: cell A6 selected, and the
# frozen-pane view scrolled back so column A is visible (no topLeftCell
# override left over from before the edit).
$ws.Range("A6").Select()
